# "code status as of green light moment / minor modifications"
#
# Re-labels a handful of coded "discussion" categories on the
# "answers" sheet (several near-duplicate categories get collapsed
# into "minimal" / "full discussion probably deleted"), marks a few
# rows as also coded "x" under evaluation/reasoning, and updates the
# remembered cell-selection on both sheets.

$wb = $excel.ActiveWorkbook
$wsAnswers = $wb.Worksheets.Item("answers")
$wsQuestions = $wb.Worksheets.Item("questions")

# --- answers sheet -------------------------------------------------

# New "x" marks added in the evaluation (C) / reasoning (E) columns.
$wsAnswers.Range("C23").Value = "x"
$wsAnswers.Range("E26").Value = "x"
$wsAnswers.Range("C29").Value = "x"
$wsAnswers.Range("E41").Value = "x"
$wsAnswers.Range("E45").Value = "x"
$wsAnswers.Range("E49").Value = "x"

# Several near-duplicate "no discussion ..." categories were merged
# away; re-point those rows at the surviving / new category labels.
$wsAnswers.Range("H23").Value = "minimal"
$wsAnswers.Range("H26").Value = "minimal"
$wsAnswers.Range("H29").Value = "full discussion probably deleted"
$wsAnswers.Range("H41").Value = "minimal"
$wsAnswers.Range("H45").Value = "minimal"
$wsAnswers.Range("H49").Value = "minimal"

# --- window / selection state ---------------------------------------
# The questions sheet was visited (scrolled/selected) and then focus
# returned to the answers sheet, which stays the active tab.

$wsQuestions.Activate()
$wsQuestions.Range("H9").Select()

$wsAnswers.Activate()
$wsAnswers.Range("H31").Select()
